# Update "Generate Report for Handback" timestamps.
# - Overview!G2  and de-de!H2 share the "Latest HO Xliff Generate Date" /
#   "Latest Handoff Datetime" value for the 0434c169... file; bump it.
# - zh-cn!H2 (Latest Handoff Datetime) and zh-cn!K2 (Latest Handback DateTime)
#   for the same file; bump them.
# - de-de!K2 (Latest Handback DateTime) for the same file; bump it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-19 01:02:17"
$wsDeDe.Range("H2").Value = "2016-08-19 01:02:17"

$wsZhCn.Range("H2").Value = "2016-08-19 01:02:12"
$wsZhCn.Range("K2").Value = "2016-08-19 01:02:30"

$wsDeDe.Range("K2").Value = "2016-08-19 01:02:37"
